# Rename the existing "Sheet1" to "codes", add a new "rooms" sheet right
# after it, and populate "rooms" with the course/room lookup table.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---------------------------------------
$codes = $wb.Worksheets.Item(1)
$codes.Name = "codes"

# --- Add the new "rooms" sheet, placed right after "codes" -----------
$rooms = $wb.Worksheets.Add($null, $codes)
$rooms.Name = "rooms"

# --- Fill in the room assignments -------------------------------------
$data = @(
    @("Preparing for IT", "T-3010"),
    @("Information Technology Concepts", "T-3030"),
    @("Information Technology Skills", "T-3030"),
    @("Programming Logic and Design", "T-3010"),
    @("Microsoft Windows Operating Systems", "T-3050"),
    @("Data Communications", "T-3080"),
    @("Microcomputer Systems Maintenance", "T-3020")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rooms.Cells.Item($row, 1).Value = $data[$i][0]
    $rooms.Cells.Item($row, 2).Value = $data[$i][1]
}

$rooms.Columns.Item(1).ColumnWidth = 42.6

# --- Restore "codes" as the active/selected tab -----------------------
$rooms.Range("A21").Select() | Out-Null
$codes.Select() | Out-Null
